# Daily auto push: a new row is inserted at row 836, dated 2026/02/17
# (continuing the date sequence already present at rows 834-835), which
# pushes the previous rows 836-877 ("2026/12/29" .. "2027/01/05") down
# to rows 837-878.
#
# Column A stores the date as literal text (not a real Excel date
# serial), exactly like every other date cell in the sheet. Typing a
# date-shaped string via .Value would make Excel auto-convert it to a
# date number, so instead the new A836 cell is seeded by copying the
# already-text "2026/02/17" cell sitting right above it (A835) - this
# carries over both the text value and its plain "General" formatting
# without introducing any new cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 836..877 down to 837..878, leaving a blank row 836.
$ws.Rows.Item(836).Insert()

# Seed A836 with the literal text "2026/02/17" by copying the identical
# text cell from A835 (avoids Excel's automatic text->date coercion).
$ws.Cells.Item(835, 1).Copy($ws.Cells.Item(836, 1))
$excel.CutCopyMode = $false

$ws.Cells.Item(836, 2).Value = "火"
$ws.Cells.Item(836, 3).Value = 15
$ws.Cells.Item(836, 4).Value = 201
